# Update meanEMG leg max ROM values (Hjemme passive) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (ROM indices) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (Subj 1 / CON) updated passive force values for columns B:E
$ws.Range("B2").Value = 826.86838530406374
$ws.Range("C2").Value = 461.50722525496786
$ws.Range("D2").Value = 975.62008339101214
$ws.Range("E2").Value = 417.24887536970186

# Row 3 (Subj 2 / STR) updated passive force values for columns B:E
$ws.Range("B3").Value = 770.17302114195684
$ws.Range("C3").Value = 444.87390981478114
$ws.Range("D3").Value = 1377.2977032292674
$ws.Range("E3").Value = 641.6741646289305

# Match the updated selection shown in the saved workbook
$ws.Range("B1:E3").Select()
